# Chapter 3 edit:
#   1. Rewrite the Introduction's purpose paragraph ("Introduce your work
#      in this chapter, what is going to be covered and how.") as a new,
#      longer sentence split across seven runs, dropping the old
#      mid-sentence _GoBack bookmark and gramStart/gramEnd proofing marks
#      along with it.
#   2. Move the _GoBack bookmark down to the very start of the next
#      paragraph ("Current System."), and drop that paragraph's stray
#      leading-space run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Find the old sentence and swap it for the new one.
# ---------------------------------------------------------------------
$oldSentence = "Introduce your work in this chapter , what is going to be covered and how."

$target = $d.Content
$found = $target.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the original introduction sentence"
}

$paraStart = $target.Start

$enDash = [char]0x2013
$segments = @(
    "This chapter ",
    "entails ",
    ("an analysis of the resources " + $enDash + " software, hardware and others " + $enDash + " "),
    "that are needed ",
    "for the proposed ",
    "system to ",
    "function or work properly."
)
$newSentence = [string]::Join("", $segments)

# Replacing the whole matched range (which spans the old runs, the
# gramStart/gramEnd proofErr marks and the old _GoBack bookmark) in one
# shot removes all of those marks and leaves a single run carrying the
# paragraph's Times New Roman / sz24 formatting.
$whole = $d.Range($target.Start, $target.End)
$whole.Text = $newSentence

# ---------------------------------------------------------------------
# 2) Re-split that single run back into seven runs - one per phrase -
#    by dropping a bookmark on each internal seam and deleting it again.
#    Word always breaks a run at a bookmark position, and the break is
#    retained even after the bookmark itself goes away.
# ---------------------------------------------------------------------
$pos = $paraStart
$seamCount = 0
for ($i = 0; $i -lt ($segments.Length - 1); $i++) {
    $pos = $pos + $segments[$i].Length
    $seamCount = $seamCount + 1
    $seam = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TempSplitSeam$seamCount", $seam)
}
for ($j = 1; $j -le $seamCount; $j++) {
    $d.Bookmarks("TempSplitSeam$j").Delete()
}

# ---------------------------------------------------------------------
# 3) The heading paragraph right after ("Current System.") loses its
#    leading standalone space run, and gains a fresh, collapsed _GoBack
#    bookmark planted at its very start (ahead of the gramStart mark).
# ---------------------------------------------------------------------
$headingPara = $target.Paragraphs(1).Next()
$headingStart = $headingPara.Range.Start

$leadingSpace = $d.Range($headingStart, $headingStart + 1)
if ($leadingSpace.Text -ne " ") {
    throw "Unexpected content before 'Current System.' heading"
}
$leadingSpace.Text = ""

$goBackSpot = $d.Range($headingStart, $headingStart)
$d.Bookmarks.Add("_GoBack", $goBackSpot)
